$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.3157894736842105
$ws.Range("C2").Value = 0.3157894736842105
$ws.Range("P2").Value = 0.1578947368421053
$ws.Range("S2").Value = 0.2105263157894737
$ws.Range("P3").Value = 0.6666666666666666
$ws.Range("S3").Value = 0.3333333333333333
$ws.Range("F6").Value = 0.2941176470588235
$ws.Range("J6").Value = 0.1176470588235294
$ws.Range("Q6").Value = 0.1176470588235294
$ws.Range("S6").Value = 0.4705882352941176
$ws.Range("Q7").Value = 0.1538461538461539
$ws.Range("S7").Value = 0.8461538461538461
$ws.Range("B8").Value = 0.08163265306122448
$ws.Range("F8").Value = 0.1020408163265306
$ws.Range("J8").Value = 0.04081632653061224
$ws.Range("O8").Value = 0.02040816326530612
$ws.Range("Q8").Value = 0.1020408163265306
$ws.Range("R8").Value = 0.1428571428571428
$ws.Range("S8").Value = 0.5102040816326531
$ws.Range("F9").Value = 0.1428571428571428
$ws.Range("Q9").Value = 0.07142857142857142
$ws.Range("R9").Value = 0.2142857142857143
$ws.Range("S9").Value = 0.5714285714285714
$ws.Range("B10").Value = 0.1216216216216216
$ws.Range("F10").Value = 0.01351351351351351
$ws.Range("J10").Value = 0.06756756756756757
$ws.Range("O10").Value = 0.01351351351351351
$ws.Range("Q10").Value = 0.0945945945945946
$ws.Range("R10").Value = 0.08108108108108109
$ws.Range("S10").Value = 0.6081081081081081
$ws.Range("G11").Value = 0.1851851851851852
$ws.Range("J11").Value = 0.03703703703703703
$ws.Range("K11").Value = 0.2592592592592592
$ws.Range("L11").Value = 0.4074074074074074
$ws.Range("S11").Value = 0.1111111111111111
$ws.Range("G12").Value = 0.5454545454545454
$ws.Range("J12").Value = 0.3636363636363636
$ws.Range("S12").Value = 0.09090909090909091
$ws.Range("G13").Value = 0.5
$ws.Range("J13").Value = 0.25
$ws.Range("S13").Value = 0.25
$ws.Range("S14").Value = 1
$ws.Range("F15").Value = 0.05882352941176471
$ws.Range("H15").Value = 0.1176470588235294
$ws.Range("I15").Value = 0.1764705882352941
$ws.Range("J15").Value = 0.4117647058823529
$ws.Range("M15").Value = 0.05882352941176471
$ws.Range("S15").Value = 0.1764705882352941
$ws.Range("H16").Value = 0.1666666666666667
$ws.Range("I16").Value = 0.3333333333333333
$ws.Range("J16").Value = 0.3333333333333333
$ws.Range("S16").Value = 0.1666666666666667
$ws.Range("H17").Value = 0.1764705882352941
$ws.Range("J17").Value = 0.2352941176470588
$ws.Range("K17").Value = 0.1176470588235294
$ws.Range("O17").Value = 0.05882352941176471
$ws.Range("S17").Value = 0.4117647058823529
$ws.Range("H18").Value = 0.2941176470588235
$ws.Range("I18").Value = 0.2352941176470588
$ws.Range("J18").Value = 0.2941176470588235
$ws.Range("K18").Value = 0.05882352941176471
$ws.Range("O18").Value = 0.05882352941176471
$ws.Range("S18").Value = 0.05882352941176471
$ws.Range("F19").Value = 0.01418439716312057
$ws.Range("H19").Value = 0.2836879432624114
$ws.Range("I19").Value = 0.04964539007092199
$ws.Range("J19").Value = 0.3049645390070922
$ws.Range("K19").Value = 0.1134751773049645
$ws.Range("M19").Value = 0.01418439716312057
$ws.Range("O19").Value = 0.07092198581560284
$ws.Range("S19").Value = 0.148936170212766
